$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Instructions")
$ws1.Range("A1").Value = "test"
